# Adding Social related Classes and Functions
# Adds a new "SA_Reports_Post" worksheet (after the existing CA_FIL sheet)
# with vendor / from-date / to-date report criteria, matching the
# FilterCriteria-style layout used by the other sheets in this workbook.

$wb = $excel.ActiveWorkbook

# Add the new worksheet as the LAST tab (after the current last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "SA_Reports_Post"

# Row 1 - headers / section titles
$ws.Range("A1").Value = "Vendor"
$ws.Range("B1").Value = "Zoom"
$ws.Range("C1").Value = "FromDate"
$ws.Range("F1").Value = "ToDate"

# Row 2 - column sub-headers (under FromDate / ToDate)
$ws.Range("C2").Value = "Day_DD"
$ws.Range("D2").Value = "Month_MMM"
$ws.Range("E2").Value = "Year_YYYY"
$ws.Range("F2").Value = "Day_DD"
$ws.Range("G2").Value = "Month_MMM"
$ws.Range("H2").Value = "Year_YYYY"

# Row 3 - sample scenario values
$ws.Range("A3").Value = "Instagram"
$ws.Range("B3").Value = "6m"
$ws.Range("D3").Value = "null"
$ws.Range("G3").Value = "null"

# Center the FromDate / ToDate section headers, then merge each across
# their three columns (C1:E1 and F1:H1).
$ws.Range("C1:H1").HorizontalAlignment = -4108
$ws.Range("C1:E1").Merge()
$ws.Range("F1:H1").Merge()

# Column widths to roughly match the rest of the workbook's report sheets.
$ws.Columns.Item(1).ColumnWidth = 9.7109375
$ws.Columns.Item(2).ColumnWidth = 6
$ws.Columns.Item(4).ColumnWidth = 13.140625
$ws.Columns.Item(5).ColumnWidth = 9.85546875
$ws.Columns.Item(6).ColumnWidth = 7.85546875
$ws.Columns.Item(7).ColumnWidth = 13.140625
$ws.Columns.Item(8).ColumnWidth = 9.42578125

# New sheet becomes the active / selected tab.
$ws.Select()
$ws.Range("G3").Select()
